$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 corresponds to diode D2 ("Schottky Barrier Diode").
# Update its Footprint and LCSC Part # to reflect the new diode choice.
$ws.Range("C5").Value = "SOD-523"
$ws.Range("D5").Value = "C345957"

# Row 6 corresponds to the 4-Pin Connector (J1,J2); its LCSC Part # changes too.
$ws.Range("D6").Value = "C145956"

# Update the active selection to match the saved view state (entire row 6 selected,
# with A6 as the active cell).
$ws.Range("A6:XFD6").Select()
